$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jalon 2")

# Update D4 (6 -> 8) and D5 (0 -> 2); G column formulas (F-D) recalc automatically
$ws.Range("D4").Value = 8
$ws.Range("D5").Value = 2

# Add new row 9 data: index 8 and new task "CSS responsive, SASS"
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "CSS responsive, SASS"

# Update the selected/active cell in the sheet view from B7 to D7
$ws.Range("D7").Select()
